$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the "RBI (India)" scenario cell to the new scenario text and give it
# a left/top aligned look (new cell style).
$cell = $wsInput.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

# Make ProductLoanInput the active sheet/tab and leave the selection on the
# cell that was just edited.
$wsInput.Activate()
$wsInput.Range("B17").Select()
